# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps that get refreshed each time the
# handback status report is (re)generated.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
# "Latest HO Xliff Generate Date" for d452f1dd-74d8-4c8f-972c-7be67665c439.md
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-21 20:56:56"

# --- zh-cn sheet -----------------------------------------------------------
# "Correspond Handoff Datetime" / "Correspond Handback DateTime" for the
# d452f1dd-74d8-4c8f-972c-7be67665c439 row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-21 20:56:52"
$wsZhCn.Range("K4").Value = "2016-08-21 20:57:13"

# --- de-de sheet -----------------------------------------------------------
# "Correspond Handoff Datetime" (mirrors the Overview value above) and
# "Correspond Handback DateTime" for the same row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-08-21 20:56:56"
$wsDeDe.Range("K4").Value = "2016-08-21 20:57:20"
